$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date range update: 17-10-2025 -> 18-10-2025 ---
$ws.Range("I1").Value = "18-10-2025 00:00:00"
$ws.Range("K1").Value = "18-10-2025 00:00:00"

# --- Row group [161, 162, 163]: rotate B/C/D/E/F/G among rows 161, 162, 163 ---
$ws.Range("B161").Value = 53925
$ws.Range("C161").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D161").Value = 66.44
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44
$ws.Range("B162").Value = 64350
$ws.Range("C162").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D162").Value = 66.44
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 79
$ws.Range("G162").Value = 5248.76
$ws.Range("B163").Value = 57756
$ws.Range("C163").Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Range("D163").Value = 66.44
$ws.Range("E163").Value = 79.37
$ws.Range("F163").Value = -100
$ws.Range("G163").Value = -6644

# --- Row group [183, 184]: rotate B/C/D/E/F/G among rows 183, 184 ---
$ws.Range("B183").Value = 57552
$ws.Range("D183").Value = 120.69
$ws.Range("E183").Value = 136.86
$ws.Range("F183").Value = -5
$ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329
$ws.Range("D184").Value = 120.69
$ws.Range("E184").Value = 128.32
$ws.Range("F184").Value = 6
$ws.Range("G184").Value = 724.14

# --- Row group [264, 265]: rotate B/C/D/E/F/G among rows 264, 265 ---
$ws.Range("B264").Value = 48719
$ws.Range("D264").Value = 295.75
$ws.Range("E264").Value = 353.35
$ws.Range("F264").Value = -81
$ws.Range("G264").Value = -23955.75
$ws.Range("B265").Value = 64979
$ws.Range("D265").Value = 295.75
$ws.Range("E265").Value = 314.41
$ws.Range("F265").Value = 62
$ws.Range("G265").Value = 18336.5

# --- Row group [279, 280]: rotate B/C/D/E/F/G among rows 279, 280 ---
$ws.Range("B279").Value = 48706
$ws.Range("D279").Value = 33.3
$ws.Range("E279").Value = 39.8
$ws.Range("F279").Value = -144
$ws.Range("G279").Value = -4795.2
$ws.Range("B280").Value = 64973
$ws.Range("D280").Value = 33.3
$ws.Range("E280").Value = 35.4
$ws.Range("F280").Value = 146
$ws.Range("G280").Value = 4861.8

# --- Row group [313, 314]: rotate B/C/D/E/F/G among rows 313, 314 ---
$ws.Range("B313").Value = 57854
$ws.Range("D313").Value = 305.84
$ws.Range("E313").Value = 325.16
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.68
$ws.Range("B314").Value = 62997
$ws.Range("D314").Value = 305.84
$ws.Range("E314").Value = 325.16
$ws.Range("F314").Value = 0
$ws.Range("G314").Value = 0

# --- Row group [317, 318]: rotate B/C/D/E/F/G among rows 317, 318 ---
$ws.Range("B317").Value = 57077
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08
$ws.Range("B318").Value = 63565
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6

# --- Row group [372, 373]: rotate B/C/D/E/F/G among rows 372, 373 ---
$ws.Range("B372").Value = 57885
$ws.Range("D372").Value = 52.13
$ws.Range("E372").Value = 62.28
$ws.Range("F372").Value = 4
$ws.Range("G372").Value = 208.52
$ws.Range("B373").Value = 63652
$ws.Range("D373").Value = 52.13
$ws.Range("E373").Value = 55.42
$ws.Range("F373").Value = 233
$ws.Range("G373").Value = 12146.29

# --- Row group [400, 401]: rotate B/C/D/E/F/G among rows 400, 401 ---
$ws.Range("B400").Value = 62933
$ws.Range("E400").Value = 70.65
$ws.Range("F400").Value = 146
$ws.Range("G400").Value = 8632.98
$ws.Range("B401").Value = 57835
$ws.Range("E401").Value = 70.65
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13

# --- Row group [431, 432]: rotate B/C/D/E/F/G among rows 431, 432 ---
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F431").Value = 4
$ws.Range("G431").Value = 237.88
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47

# --- Row group [536, 537]: rotate B/C/D/E/F/G among rows 536, 537 ---
$ws.Range("B536").Value = 58047
$ws.Range("D536").Value = 105.54
$ws.Range("E536").Value = 126.1
$ws.Range("F536").Value = 52
$ws.Range("G536").Value = 5488.08
$ws.Range("B537").Value = 47097
$ws.Range("D537").Value = 112.28
$ws.Range("E537").Value = 134.16
$ws.Range("F537").Value = 15
$ws.Range("G537").Value = 1684.2

# --- Row group [583, 584]: rotate B/C/D/E/F/G among rows 583, 584 ---
$ws.Range("B583").Value = 53263
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29
$ws.Range("B584").Value = 65066
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 278
$ws.Range("G584").Value = 3561.18

# --- Row group [599, 600]: rotate B/C/D/E/F/G among rows 599, 600 ---
$ws.Range("B599").Value = 45709
$ws.Range("E599").Value = 15.69
$ws.Range("F599").Value = -300
$ws.Range("G599").Value = -3945
$ws.Range("B600").Value = 64925
$ws.Range("E600").Value = 13.97
$ws.Range("F600").Value = 281
$ws.Range("G600").Value = 3695.15

# --- Row group [601, 602]: rotate B/C/D/E/F/G among rows 601, 602 ---
$ws.Range("B601").Value = 64919
$ws.Range("E601").Value = 27.97
$ws.Range("F601").Value = 209
$ws.Range("G601").Value = 5496.7
$ws.Range("B602").Value = 45702
$ws.Range("E602").Value = 31.43
$ws.Range("F602").Value = -215
$ws.Range("G602").Value = -5654.5

# --- Row group [709, 710]: rotate B/C/D/E/F/G among rows 709, 710 ---
$ws.Range("B709").Value = 64833
$ws.Range("D709").Value = 32.83
$ws.Range("E709").Value = 34.9
$ws.Range("F709").Value = 99
$ws.Range("G709").Value = 3250.17
$ws.Range("B710").Value = 60025
$ws.Range("D710").Value = 32.83
$ws.Range("E710").Value = 37.22
$ws.Range("F710").Value = -98
$ws.Range("G710").Value = -3217.34

# --- Row group [715, 716]: rotate B/C/D/E/F/G among rows 715, 716 ---
$ws.Range("B715").Value = 60031
$ws.Range("D715").Value = 98.5
$ws.Range("E715").Value = 111.69
$ws.Range("F715").Value = -5
$ws.Range("G715").Value = -492.5
$ws.Range("B716").Value = 64836
$ws.Range("D716").Value = 98.5
$ws.Range("E716").Value = 104.71
$ws.Range("F716").Value = 7
$ws.Range("G716").Value = 689.5
